$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 177
$ws.Cells.Item(177, 1).Value2() = 177
$ws.Cells.Item(177, 2).Value2() = "CESGRANRIO"
$ws.Cells.Item(177, 3).Value2() = "TRANSPETRO"
$ws.Cells.Item(177, 4).Value2() = "'2018"
$ws.Cells.Item(177, 5).Value2() = "Os sinais de pontuação contribuem para a construção dos sentidos dos textos. No fragmento “Escriturei-me; deram-me um papel que... mas para que o estou a fatigar com isso? Deixe-me ficar com as minhas amofinações”, as reticências são usadas para demarcar a"
$ws.Cells.Item(177, 6).Value2() = "Português"
$ws.Cells.Item(177, 7).Value2() = "Reticências"
$ws.Cells.Item(177, 8).Value2() = "Médio"
$ws.Cells.Item(177, 9).Value2() = "ME"
$ws.Cells.Item(177, 10).Value2() = "interrupção de uma ideia."
$ws.Cells.Item(177, 11).Value2() = "insinuação de uma ameaça."
$ws.Cells.Item(177, 12).Value2() = "hesitação comum na oralidade."
$ws.Cells.Item(177, 13).Value2() = "continuidade de uma ação ou fato."
$ws.Cells.Item(177, 14).Value2() = "omissão proposital de algo que se devia dizer."
$ws.Cells.Item(177, 15).Value2() = "A"
$ws.Cells.Item(177, 16).Value2() = 0
$ws.Cells.Item(177, 17).Value2() = 0

# Row 178
$ws.Cells.Item(178, 1).Value2() = 178
$ws.Cells.Item(178, 2).Value2() = "CESGRANRIO"
$ws.Cells.Item(178, 3).Value2() = "PETROBRAS"
$ws.Cells.Item(178, 4).Value2() = "'2014"
$ws.Cells.Item(178, 5).Value2() = "As reticências utilizadas pelo autor no trecho “desabotoava a blusa até o estômago, enfiava a mão dentro dela e puxava para fora um seio lindo, liso, branco, aquele mamilo atrevido... E nós, meninos, de boca aberta...” assinalam uma determinada sensação.`nO trecho em que semelhante sensação se verifica é:"
$ws.Cells.Item(178, 6).Value2() = "Português"
$ws.Cells.Item(178, 7).Value2() = "Reticências"
$ws.Cells.Item(178, 8).Value2() = "Médio"
$ws.Cells.Item(178, 9).Value2() = "ME"
$ws.Cells.Item(178, 10).Value2() = "“Se estou com fome e gosto de queijo, eu como queijo...”"
$ws.Cells.Item(178, 11).Value2() = "“Procuro outra coisa de que goste: banana, pão com manteiga, chocolate...”"
$ws.Cells.Item(178, 12).Value2() = "“Enquanto varria e limpava, sofria ouvindo o pianista tocando uma música horrível: Bach, Brahms, Debussy...”"
$ws.Cells.Item(178, 13).Value2() = "“Mas eu comprava um mata-fome e ia para casa comendo o mata-fome bem devagarzinho...”"
$ws.Cells.Item(178, 14).Value2() = "“Ridendo dicere severum: rindo, dizer as coisas sérias...”"
$ws.Cells.Item(178, 15).Value2() = "D"
$ws.Cells.Item(178, 16).Value2() = 0
$ws.Cells.Item(178, 17).Value2() = 0
